# Refresh cached market-price-derived profit figures (scheduled runner update)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 457.35294
$ws.Cells.Item(2, 9).Value = 157.2
$ws.Cells.Item(2, 10).Value = 886.1429000000001
$ws.Cells.Item(2, 11).Value = 157.2
$ws.Cells.Item(2, 12).Value = 886.1429000000001
$ws.Cells.Item(2, 13).Value = -44.19999999999999
$ws.Cells.Item(2, 14).Value = -1112.1429

$ws.Cells.Item(31, 8).Value = 7314.3
$ws.Cells.Item(31, 9).Value = 383.33334
$ws.Cells.Item(31, 10).Value = 10284.714
$ws.Cells.Item(31, 11).Value = 1150.00002
$ws.Cells.Item(31, 12).Value = 30854.142
$ws.Cells.Item(31, 13).Value = -920.0000199999999
$ws.Cells.Item(31, 14).Value = -31314.142

$ws.Cells.Item(51, 8).Value = 5087.095
$ws.Cells.Item(51, 10).Value = 4691.1113
$ws.Cells.Item(51, 12).Value = 4691.1113
$ws.Cells.Item(51, 14).Value = -5659.1113

$ws.Cells.Item(138, 8).Value = 2541.4912
$ws.Cells.Item(138, 9).Value = 1377.7188
$ws.Cells.Item(138, 11).Value = 4133.1564
$ws.Cells.Item(138, 13).Value = 1006.8436

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2906.422
$ws.Cells.Item(32, 10).Value = 2000
$ws.Cells.Item(32, 12).Value = 2000
$ws.Cells.Item(32, 14).Value = -2574

$ws.Cells.Item(61, 8).Value = 50003020
$ws.Cells.Item(61, 9).Value = 55558410
$ws.Cells.Item(61, 11).Value = 55558410
$ws.Cells.Item(61, 13).Value = -55558198

$ws.Cells.Item(74, 8).Value = 45462716
$ws.Cells.Item(74, 9).Value = 52639200
$ws.Cells.Item(74, 11).Value = 52639200
$ws.Cells.Item(74, 13).Value = -52638326

$ws.Cells.Item(77, 8).Value = 45462716
$ws.Cells.Item(77, 9).Value = 52639200
$ws.Cells.Item(77, 11).Value = 263196000
$ws.Cells.Item(77, 13).Value = -263191632

$ws.Cells.Item(122, 8).Value = 2624.875
$ws.Cells.Item(122, 9).Value = 2535.5715
$ws.Cells.Item(122, 11).Value = 7606.7145
$ws.Cells.Item(122, 13).Value = -5156.7145

$ws.Cells.Item(132, 8).Value = 3229023.2
$ws.Cells.Item(132, 9).Value = 4003252.2
$ws.Cells.Item(132, 11).Value = 12009756.6
$ws.Cells.Item(132, 13).Value = -12007226.6

$ws.Cells.Item(136, 8).Value = 50003020
$ws.Cells.Item(136, 9).Value = 55558410
$ws.Cells.Item(136, 11).Value = 166675230
$ws.Cells.Item(136, 13).Value = -166672680

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 11803.053
$ws.Cells.Item(94, 9).Value = 11686.611
$ws.Cells.Item(94, 11).Value = 11686.611
$ws.Cells.Item(94, 13).Value = -11235.611

$ws.Cells.Item(134, 8).Value = 85834650
$ws.Cells.Item(134, 9).Value = 85834650
$ws.Cells.Item(134, 11).Value = 257503950
$ws.Cells.Item(134, 13).Value = -257501415

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 16671041
$ws.Cells.Item(58, 9).Value = 21744342
$ws.Cells.Item(58, 11).Value = 21744342
$ws.Cells.Item(58, 13).Value = -21744139

$ws.Cells.Item(100, 8).Value = 0
$ws.Cells.Item(100, 10).Value = 0
$ws.Cells.Item(100, 12).Value = 0

$ws.Cells.Item(105, 8).Value = 2501625.5
$ws.Cells.Item(105, 10).Value = 2422
$ws.Cells.Item(105, 12).Value = 2422
$ws.Cells.Item(105, 14).Value = -5916

$ws.Cells.Item(109, 8).Value = 51122.5
$ws.Cells.Item(109, 10).Value = 51122.5
$ws.Cells.Item(109, 12).Value = 51122.5
$ws.Cells.Item(109, 14).Value = -53202.5

$ws.Cells.Item(132, 8).Value = 71431660
$ws.Cells.Item(132, 9).Value = 71431660
$ws.Cells.Item(132, 11).Value = 214294980
$ws.Cells.Item(132, 13).Value = -214292450

$ws.Cells.Item(136, 8).Value = 16671041
$ws.Cells.Item(136, 9).Value = 21744342
$ws.Cells.Item(136, 11).Value = 65233026
$ws.Cells.Item(136, 13).Value = -65230476

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 144535.42
$ws.Cells.Item(5, 10).Value = 2190
$ws.Cells.Item(5, 12).Value = 6570
$ws.Cells.Item(5, 14).Value = -6794

$ws.Cells.Item(37, 8).Value = 149976
$ws.Cells.Item(37, 10).Value = 149976
$ws.Cells.Item(37, 12).Value = 449928
$ws.Cells.Item(37, 14).Value = -450152

$ws.Cells.Item(70, 8).Value = 9330.416999999999
$ws.Cells.Item(70, 9).Value = 3995.625
$ws.Cells.Item(70, 11).Value = 11986.875
$ws.Cells.Item(70, 13).Value = -11671.875

$ws.Cells.Item(73, 8).Value = 9330.416999999999
$ws.Cells.Item(73, 9).Value = 3995.625
$ws.Cells.Item(73, 11).Value = 11986.875
$ws.Cells.Item(73, 13).Value = -10894.875

$ws.Cells.Item(121, 8).Value = 1371731.1
$ws.Cells.Item(121, 9).Value = 339999.66
$ws.Cells.Item(121, 10).Value = 1887596.9
$ws.Cells.Item(121, 11).Value = 1019998.98
$ws.Cells.Item(121, 12).Value = 5662790.699999999
$ws.Cells.Item(121, 13).Value = -1018688.98
$ws.Cells.Item(121, 14).Value = -5665410.699999999

$ws.Cells.Item(132, 8).Value = 2099.5
$ws.Cells.Item(132, 9).Value = 2099.5
$ws.Cells.Item(132, 11).Value = 18895.5
$ws.Cells.Item(132, 13).Value = -16365.5

$ws.Cells.Item(135, 8).Value = 144535.42
$ws.Cells.Item(135, 10).Value = 2190
$ws.Cells.Item(135, 12).Value = 19710
$ws.Cells.Item(135, 14).Value = -24780

$ws.Cells.Item(137, 8).Value = 50001500
$ws.Cells.Item(137, 9).Value = 50001500
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 150004500
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 13).Value = -149999400

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(21, 8).Value = 2701
$ws.Cells.Item(21, 9).Value = 2701
$ws.Cells.Item(21, 11).Value = 2701
$ws.Cells.Item(21, 13).Value = -2528

$ws.Cells.Item(30, 8).Value = 2701
$ws.Cells.Item(30, 9).Value = 2701
$ws.Cells.Item(30, 11).Value = 2701
$ws.Cells.Item(30, 13).Value = -2596

$ws.Cells.Item(63, 8).Value = 39000
$ws.Cells.Item(63, 9).Value = 39000
$ws.Cells.Item(63, 11).Value = 39000
$ws.Cells.Item(63, 13).Value = -38314

$ws.Cells.Item(66, 8).Value = 39000
$ws.Cells.Item(66, 9).Value = 39000
$ws.Cells.Item(66, 11).Value = 117000
$ws.Cells.Item(66, 13).Value = -113568

$ws.Cells.Item(80, 8).Value = 3304.182
$ws.Cells.Item(80, 9).Value = 2356
$ws.Cells.Item(80, 11).Value = 2356
$ws.Cells.Item(80, 13).Value = -1358

$ws.Cells.Item(83, 8).Value = 3304.182
$ws.Cells.Item(83, 9).Value = 2356
$ws.Cells.Item(83, 11).Value = 11780
$ws.Cells.Item(83, 13).Value = -6788

$ws.Cells.Item(102, 8).Value = 7796.857
$ws.Cells.Item(102, 9).Value = 1596.3334
$ws.Cells.Item(102, 11).Value = 1596.3334
$ws.Cells.Item(102, 13).Value = 25.66660000000002

$ws.Cells.Item(132, 8).Value = 25001752
$ws.Cells.Item(132, 9).Value = 25001752
$ws.Cells.Item(132, 11).Value = 75005256
$ws.Cells.Item(132, 13).Value = -75002726

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(6, 8).Value = 59900
$ws.Cells.Item(6, 10).Value = 59900
$ws.Cells.Item(6, 12).Value = 59900
$ws.Cells.Item(6, 14).Value = -60124

$ws.Cells.Item(132, 8).Value = 16005725
$ws.Cells.Item(132, 9).Value = 19205890
$ws.Cells.Item(132, 11).Value = 57617670
$ws.Cells.Item(132, 13).Value = -57615140

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 31829.4
$ws.Cells.Item(2, 9).Value = 32287.5
$ws.Cells.Item(2, 11).Value = 32287.5
$ws.Cells.Item(2, 13).Value = -32175.5

$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 11).Value = 0

$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 11).Value = 0

$ws.Cells.Item(127, 8).Value = 90195
$ws.Cells.Item(127, 9).Value = 80390
$ws.Cells.Item(127, 11).Value = 80390
$ws.Cells.Item(127, 13).Value = -75430

$ws.Cells.Item(132, 8).Value = 35724956
$ws.Cells.Item(132, 9).Value = 55559716
$ws.Cells.Item(132, 11).Value = 166679148
$ws.Cells.Item(132, 13).Value = -166676618

$ws.Cells.Item(136, 8).Value = 10871225
$ws.Cells.Item(136, 9).Value = 12196936
$ws.Cells.Item(136, 11).Value = 36590808
$ws.Cells.Item(136, 13).Value = -36588258

# Clear cells with no longer meaningful cached profit value
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(100, 14).ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(137, 14).ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(64, 13).ClearContents()
$ws.Cells.Item(67, 13).ClearContents()
